$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H28").Value = 1414.9286
$ws.Range("I28").Value = 1132.3334
$ws.Range("K28").Value = 1132.3334
$ws.Range("M28").Value = -647.3334
$ws.Range("H41").Value = 1670.0714
$ws.Range("I41").Value = 153.66667
$ws.Range("J41").Value = 2807.375
$ws.Range("K41").Value = 153.66667
$ws.Range("L41").Value = 2807.375
$ws.Range("M41").Value = 286.33333
$ws.Range("N41").Value = -3687.375
$ws.Range("H51").Value = 9581.25
$ws.Range("J51").Value = 9698.299999999999
$ws.Range("L51").Value = 9698.299999999999
$ws.Range("N51").Value = -10666.3
$ws.Range("H92").Value = 467.73334
$ws.Range("I92").Value = 377.76923
$ws.Range("K92").Value = 377.76923
$ws.Range("M92").Value = 870.23077
$ws.Range("H107").Value = 1007.8333
$ws.Range("I107").Value = 1081.6364
$ws.Range("K107").Value = 1081.6364
$ws.Range("M107").Value = 838.3635999999999
$ws.Range("H131").Value = 2133.9375
$ws.Range("I131").Value = 742.93335
$ws.Range("J131").Value = 22999
$ws.Range("K131").Value = 2228.80005
$ws.Range("L131").Value = 68997
$ws.Range("M131").Value = 2811.19995
$ws.Range("N131").Value = -79077
$ws.Range("H132").Value = 6055.136
$ws.Range("I132").Value = 1227.4
$ws.Range("K132").Value = 3682.2
$ws.Range("M132").Value = -1152.2

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 51351
$ws.Range("J2").Value = 1500
$ws.Range("L2").Value = 1500
$ws.Range("N2").Value = -1726
$ws.Range("H116").Value = 51351
$ws.Range("J116").Value = 1500
$ws.Range("L116").Value = 1500
$ws.Range("N116").Value = -6088
$ws.Range("H122").Value = 8304.210999999999
$ws.Range("I122").Value = 7702.778
$ws.Range("J122").Value = 8845.5
$ws.Range("K122").Value = 23108.334
$ws.Range("L122").Value = 26536.5
$ws.Range("M122").Value = -20658.334
$ws.Range("N122").Value = -31436.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 51351
$ws.Range("J3").Value = 1500
$ws.Range("L3").Value = 1500
$ws.Range("N3").Value = -1728
$ws.Range("H99").Value = 1708.6666
$ws.Range("I99").Value = 1592.4706
$ws.Range("J99").Value = 2202.5
$ws.Range("K99").Value = 1592.4706
$ws.Range("L99").Value = 2202.5
$ws.Range("M99").Value = -94.4706000000001
$ws.Range("N99").Value = -5198.5
$ws.Range("H107").Value = 1527.3334
$ws.Range("J107").Value = 1948.8667
$ws.Range("L107").Value = 1948.8667
$ws.Range("N107").Value = -5788.8667

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2214
$ws.Range("I16").Value = 1999.75
$ws.Range("J16").Value = 2499.6667
$ws.Range("K16").Value = 1999.75
$ws.Range("L16").Value = 2499.6667
$ws.Range("M16").Value = -1712.75
$ws.Range("N16").Value = -3073.6667
$ws.Range("H113").Value = 2214
$ws.Range("I113").Value = 1999.75
$ws.Range("J113").Value = 2499.6667
$ws.Range("K113").Value = 1999.75
$ws.Range("L113").Value = 2499.6667
$ws.Range("M113").Value = 170.25
$ws.Range("N113").Value = -6839.6667
$ws.Range("H122").Value = 3987.3076
$ws.Range("I122").Value = 2735.2632
$ws.Range("K122").Value = 8205.7896
$ws.Range("M122").Value = -5755.7896

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2119.0667
$ws.Range("J114").Value = 2313
$ws.Range("L114").Value = 6939
$ws.Range("N114").Value = -13447
$ws.Range("H137").Value = 94446270
$ws.Range("J137").Value = 20002044
$ws.Range("L137").Value = 60006132
$ws.Range("N137").Value = -60016332
$ws.Range("H141").Value = 5993.091
$ws.Range("I141").Value = 2820.8333
$ws.Range("K141").Value = 8462.499899999999
$ws.Range("M141").Value = -3282.499899999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1366.2222
$ws.Range("I97").Value = 1162.125
$ws.Range("K97").Value = 1162.125
$ws.Range("M97").Value = -666.125
$ws.Range("H122").Value = 404684.75
$ws.Range("I122").Value = 911738.5600000001
$ws.Range("J122").Value = 6285.357
$ws.Range("K122").Value = 2735215.68
$ws.Range("L122").Value = 18856.071
$ws.Range("M122").Value = -2732765.68
$ws.Range("N122").Value = -23756.071
$ws.Range("H132").Value = 4067.4285
$ws.Range("I132").Value = 3192.3635
$ws.Range("J132").Value = 7276
$ws.Range("K132").Value = 9577.0905
$ws.Range("L132").Value = 21828
$ws.Range("M132").Value = -7047.0905
$ws.Range("N132").Value = -26888

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I40").Value = 3633.5
$ws.Range("K40").Value = 3633.5
$ws.Range("M40").Value = -3497.5
$ws.Range("H43").Value = 18800
$ws.Range("I43").Value = 18800
$ws.Range("K43").Value = 18800
$ws.Range("M43").Value = -18607
$ws.Range("H93").Value = 1340.2593
$ws.Range("I93").Value = 1172.8334
$ws.Range("K93").Value = 1172.8334
$ws.Range("M93").Value = 75.16660000000002
$ws.Range("H122").Value = 6142.931
$ws.Range("I122").Value = 3449.5925
$ws.Range("J122").Value = 8488.742
$ws.Range("K122").Value = 10348.7775
$ws.Range("L122").Value = 25466.226
$ws.Range("M122").Value = -7898.7775
$ws.Range("N122").Value = -30366.226
$ws.Range("H138").Value = 100429
$ws.Range("J138").Value = 100429
$ws.Range("L138").Value = 100429
$ws.Range("N138").Value = -110709

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H41").Value = 20266.111
$ws.Range("J41").Value = 20799.375
$ws.Range("L41").Value = 20799.375
$ws.Range("N41").Value = -21579.375
$ws.Range("H122").Value = 3477.5
$ws.Range("I122").Value = 2833.7144
$ws.Range("K122").Value = 8501.143199999999
$ws.Range("M122").Value = -6051.143199999999
$ws.Range("H138").Value = 91239.60000000001
$ws.Range("I138").Value = 89699
$ws.Range("J138").Value = 91624.75
$ws.Range("K138").Value = 89699
$ws.Range("L138").Value = 91624.75
$ws.Range("M138").Value = -84559
$ws.Range("N138").Value = -101904.75
$ws.Range("H140").Value = 86337.2
$ws.Range("J140").Value = 86337.2
$ws.Range("L140").Value = 86337.2
$ws.Range("N140").Value = -96697.2
